$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "54.214.75"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.260.47"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "493.97"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "127.94"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "0.0956"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "0.327"
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").Value = "4.70"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "2.660.30"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "21.91"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "54.128.82"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "2.270.67"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D19").Value = "4.10"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  +5.21%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "301.45"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "61.98"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("B24").Value = "Binance-PegBSC-USD"
$ws.Range("C24").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "2.363.40"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("D26").Value = "0.370"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.148"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "168.62"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.60"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0679"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.08"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "17.68"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.890"
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.18"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.69"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "35.76"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.40"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.371"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "127.11"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "4.77"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0891"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0484"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").Value = "0.544"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "237.52"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0203"
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.76"
$ws.Range("E51").Value = "  +0.94%  "
